$wb = $excel.ActiveWorkbook

# Sheet 1: "ورودی واگن یا بار" - update the destination (مقصد) value
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "کرمانشاه"

# Sheet 2: "خروجی" - route renamed from نورآباد leg to ملایر leg,
# gabari type becomes "نامشخص" (unknown), and the now-unreachable
# permission/size columns for that row are cleared to blank.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A5").Value = "سواریان - ملایر"
$ws2.Range("A6").Value = "ملایر - کرمانشاه"
$ws2.Range("B6").Value = "نامشخص"
$ws2.Range("C6").Value = ""
$ws2.Range("D6").Value = ""
$ws2.Range("E6").Value = ""
$ws2.Range("F6").Value = ""
